$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("İş Takip Listesi")
$ws2 = $wb.Worksheets.Item("Güncelleme")

# Helper: write a literal text value into a cell without letting Excel
# auto-convert a date-shaped string ("YYYY-MM-DD") into a real date serial.
function Set-TextCell {
    param($Sheet, $Cell, $Text)
    $range = $Sheet.Range($Cell)
    $range.NumberFormat = "@"
    $range.Value = $Text
}

# --- "İş Takip Listesi": İŞE BAŞLAMA/YER TESLİMİ (J) & İHALE BİTİŞ TARİHİ (K) ---
# Three contiguous row blocks, each shifted back by one calendar day.

for ($r = 2; $r -le 10; $r++) {
    Set-TextCell $ws1 "J$r" "2025-08-27"
    Set-TextCell $ws1 "K$r" "2026-01-27"
}

for ($r = 33; $r -le 94; $r++) {
    Set-TextCell $ws1 "J$r" "2025-08-29"
    Set-TextCell $ws1 "K$r" "2026-01-29"
}

for ($r = 95; $r -le 122; $r++) {
    Set-TextCell $ws1 "J$r" "2024-06-27"
    Set-TextCell $ws1 "K$r" "2025-08-21"
}

# Row 121 status moved on from değerlendirme to kroki/tutanak kontrolü.
$ws1.Range("L121").Value = "KROKİ/TUTANAK KONTROLÜ"

# --- "Güncelleme": scattered İ/J/N/P tracking dates, each shifted back one day ---
$guncellemeEdits = @(
    @{ Cell = "J2"; Value = "2024-10-03" }
    @{ Cell = "N2"; Value = "2025-06-06" }
    @{ Cell = "P2"; Value = "2025-08-23" }
    @{ Cell = "J3"; Value = "2025-01-04" }
    @{ Cell = "N3"; Value = "2025-09-24" }
    @{ Cell = "J4"; Value = "2024-11-08" }
    @{ Cell = "N4"; Value = "2025-05-01" }
    @{ Cell = "P4"; Value = "2025-07-26" }
    @{ Cell = "I5"; Value = "2025-05-04" }
    @{ Cell = "J6"; Value = "2025-12-15" }
    @{ Cell = "N6"; Value = "2025-09-03" }
    @{ Cell = "I7"; Value = "2025-01-04" }
    @{ Cell = "J7"; Value = "2025-01-04" }
    @{ Cell = "J8"; Value = "2024-12-23" }
    @{ Cell = "N8"; Value = "2025-05-24" }
    @{ Cell = "P8"; Value = "2025-06-26" }
    @{ Cell = "I9"; Value = "2025-08-20" }
    @{ Cell = "J9"; Value = "2025-02-05" }
    @{ Cell = "J10"; Value = "2024-12-04" }
    @{ Cell = "N10"; Value = "2025-09-14" }
    @{ Cell = "I11"; Value = "2025-06-11" }
    @{ Cell = "J11"; Value = "2025-01-17" }
    @{ Cell = "N11"; Value = "2025-10-04" }
    @{ Cell = "J12"; Value = "2024-12-15" }
    @{ Cell = "N12"; Value = "2025-08-24" }
    @{ Cell = "J13"; Value = "2025-02-12" }
    @{ Cell = "J14"; Value = "2025-12-11" }
    @{ Cell = "J15"; Value = "2025-03-03" }
    @{ Cell = "N15"; Value = "2025-09-21" }
    @{ Cell = "J16"; Value = "2024-10-30" }
    @{ Cell = "N16"; Value = "2025-04-09" }
    @{ Cell = "P16"; Value = "2025-06-26" }
    @{ Cell = "J17"; Value = "2024-11-15" }
    @{ Cell = "J18"; Value = "2025-04-24" }
    @{ Cell = "I19"; Value = "2025-06-12" }
    @{ Cell = "J19"; Value = "2025-03-03" }
    @{ Cell = "N19"; Value = "2025-10-11" }
    @{ Cell = "J20"; Value = "2025-02-12" }
    @{ Cell = "J21"; Value = "2024-12-06" }
    @{ Cell = "J22"; Value = "2024-12-06" }
    @{ Cell = "J23"; Value = "2025-02-13" }
    @{ Cell = "I24"; Value = "2025-08-10" }
    @{ Cell = "J25"; Value = "2025-01-08" }
    @{ Cell = "J27"; Value = "2025-03-31" }
    @{ Cell = "J28"; Value = "2025-01-27" }
    @{ Cell = "I29"; Value = "2025-04-18" }
    @{ Cell = "J29"; Value = "2025-02-13" }
)

foreach ($edit in $guncellemeEdits) {
    Set-TextCell $ws2 $edit.Cell $edit.Value
}
